$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting D:K to E:L
$ws.Columns("D:D").Insert()

# Copy number formatting (date / number styles) from column E into the new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the new (most recent) period figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 6705000
$ws.Range("D9").Value = 6187000
$ws.Range("D10").Value = 518000
$ws.Range("D12").Value = 20000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 2360000
$ws.Range("D15").Value = 62000
$ws.Range("D17").Value = 8961000
$ws.Range("D18").Value = -2256000
$ws.Range("D20").Value = -56000
$ws.Range("D21").Value = -2033000
$ws.Range("D22").Value = 259000
$ws.Range("D23").Value = -2571000
$ws.Range("D24").Value = 104000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -2675000
$ws.Range("D27").Value = -2691000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 56000
$ws.Range("D33").Value = -2691000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -2691000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 520000
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 1811000
$ws.Range("D44").Value = 101000
$ws.Range("D45").Value = 601000
$ws.Range("D46").Value = 3033000
$ws.Range("D47").Value = 514000
$ws.Range("D48").Value = 2067000
$ws.Range("D49").Value = 3663000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 163000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 9440000
$ws.Range("D57").Value = 595000
$ws.Range("D58").Value = 38000
$ws.Range("D59").Value = 3584000
$ws.Range("D60").Value = 4217000
$ws.Range("D61").Value = 3459000
$ws.Range("D62").Value = 711000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 8640000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -2719000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 800000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -2691000
$ws.Range("D83").Value = 279000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -71000
$ws.Range("D91").Value = -86000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2648000
$ws.Range("D96").Value = -3000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 3201000
$ws.Range("D101").Value = -45000
$ws.Range("D102").Value = 437000

# Apply restated figures to a handful of existing cells that changed value
# (not a pure shift - these numbers were revised in the source data)
$ws.Range("E8").Value = 2985000
$ws.Range("E9").Value = 2449000
$ws.Range("F9").Value = 2249000
$ws.Range("E10").Value = 536000
$ws.Range("F10").Value = 387000
$ws.Range("E12").Value = 5000
$ws.Range("E14").Value = 1000
$ws.Range("F14").Value = 132200
$ws.Range("G15").Value = "NA"
$ws.Range("H15").Value = "NA"
$ws.Range("I15").Value = "NA"
$ws.Range("J15").Value = "NA"
$ws.Range("E17").Value = 2678000
$ws.Range("F17").Value = 2498000
$ws.Range("E18").Value = 307000
$ws.Range("F18").Value = 138000
$ws.Range("E20").Value = 68000
$ws.Range("F20").Value = 64000
$ws.Range("E21").Value = 475700
$ws.Range("F21").Value = 304700
$ws.Range("E22").Value = 126000
$ws.Range("F22").Value = 120000
$ws.Range("E23").Value = 249000
$ws.Range("F23").Value = 82000
$ws.Range("E24").Value = 69000
$ws.Range("F24").Value = 42000
$ws.Range("E26").Value = 180000
$ws.Range("F26").Value = 40000
$ws.Range("E27").Value = 179000
$ws.Range("F27").Value = 34000
$ws.Range("E29").Value = 0
$ws.Range("E32").Value = -68000
$ws.Range("F32").Value = -64000
$ws.Range("E33").Value = 179000
$ws.Range("F33").Value = 34000
$ws.Range("E35").Value = 179000
$ws.Range("F35").Value = 34000
$ws.Range("E44").Value = "NA"
$ws.Range("F44").Value = "NA"
$ws.Range("G44").Value = "NA"
$ws.Range("H44").Value = "NA"
$ws.Range("I44").Value = "NA"
$ws.Range("J44").Value = "NA"
$ws.Range("E81").Value = 179000
$ws.Range("F81").Value = 34000
